# Rename the inline logo pictures that live in the document's headers and
# footers:
#   - the BTEC logo (alt text "BTec_Logo-Orange"), currently named
#     "image2.jpg", should become "image1.jpg"
#   - the two Pearson logos (alt text "...PearsonLogo.png"), currently
#     named "image1.png", should become "image2.png"

$d = $word.ActiveDocument

function Rename-LogoShapes {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $alt = $shp.AlternativeText

        if ($alt -eq "BTec_Logo-Orange") {
            $shp.Name = "image1.jpg"
        } elseif ($alt -like "*PearsonLogo.png") {
            $shp.Name = "image2.png"
        }
    }
}

foreach ($sec in $d.Sections) {
    $headers = $sec.Headers
    for ($i = 1; $i -le $headers.Count; $i++) {
        Rename-LogoShapes $headers.Item($i).Range.InlineShapes
    }

    $footers = $sec.Footers
    for ($i = 1; $i -le $footers.Count; $i++) {
        Rename-LogoShapes $footers.Item($i).Range.InlineShapes
    }
}

Rename-LogoShapes $d.InlineShapes
